# Updated cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "26.401.98"
$ws.Range("E2").Value2 = "  -1.24%  "

# Row 3
$ws.Range("D3").Value2 = "1.624.11"
$ws.Range("E3").Value2 = "  -0.80%  "

# Row 4
$ws.Range("E4").Value2 = "  +0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "212.78"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.497"
$ws.Range("E6").Value2 = "  +1.28%  "

# Row 7
$ws.Range("E7").Value2 = "  +0.22%  "

# Row 8
$ws.Range("E8").Value2 = "  -0.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.0622"
$ws.Range("E9").Value2 = "  +0.22%  "

# Row 10
$ws.Range("E10").Value2 = "  -0.53%  "

# Row 11
$ws.Range("E11").Value2 = "  +0.91%  "

# Row 12
$ws.Range("D12").Value2 = "1.849.32"
$ws.Range("E12").Value2 = "  -0.80%  "

# Row 13
$ws.Range("B13").Value2 = "Polkadot"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.13"
$ws.Range("E13").Value2 = "  +1.87%  "

# Row 14
$ws.Range("B14").Value2 = "WrappedEther"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "1.582.46"
$ws.Range("E14").Value2 = "  -3.23%  "

# Row 15
$ws.Range("E15").Value2 = "  -0.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "64.48"
$ws.Range("E16").Value2 = "  +2.06%  "

# Row 17
$ws.Range("D17").Value2 = "26.438.86"
$ws.Range("E17").Value2 = "  -0.99%  "

# Row 18
$ws.Range("E18").Value2 = "  +0.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "215.80"
$ws.Range("E19").Value2 = "  +3.47%  "

# Row 20
$ws.Range("E20").Value2 = "  +0.21%  "

# Row 21
$ws.Range("E21").Value2 = "  -0.81%  "

# Row 22
$ws.Range("E22").Value2 = "  +1.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "9.28"
$ws.Range("E23").Value2 = "  -1.15%  "

# Row 24
$ws.Range("E24").Value2 = "  +3.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "147.54"
$ws.Range("E25").Value2 = "  +1.02%  "

# Row 26
$ws.Range("E26").Value2 = "  +0.17%  "

# Row 27
$ws.Range("E27").Value2 = "  -1.03%  "

# Row 28
$ws.Range("E28").Value2 = "  +2.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "15.57"
$ws.Range("E29").Value2 = "  +1.34%  "

# Row 30
$ws.Range("E30").Value2 = "  -1.64%  "

# Row 31
$ws.Range("E31").Value2 = "  -1.33%  "

# Row 32
$ws.Range("E32").Value2 = "  +2.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "2.94"
$ws.Range("E33").Value2 = "  -0.27%  "

# Row 34
$ws.Range("E34").Value2 = "  -0.80%  "

# Row 35
$ws.Range("D35").Value2 = "1.218.14"
$ws.Range("E35").Value2 = "  +4.65%  "

# Row 36
$ws.Range("E36").Value2 = "  -1.87%  "

# Row 37
$ws.Range("E37").Value2 = "  +3.50%  "

# Row 38
$ws.Range("E38").Value2 = "  +0.18%  "

# Row 39
$ws.Range("E39").Value2 = "  -2.24%  "

# Row 40
$ws.Range("E40").Value2 = "  -0.14%  "

# Row 41
$ws.Range("E41").Value2 = "  -2.97%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.792"
$ws.Range("E42").Value2 = "  -0.44%  "

# Row 43
$ws.Range("E43").Value2 = "  -0.50%  "

# Row 44
$ws.Range("D44").Value2 = "1.762.10"
$ws.Range("E44").Value2 = "  -0.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "92.69"
$ws.Range("E45").Value2 = "  +0.14%  "

# Row 46
$ws.Range("E46").Value2 = "  +1.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "54.68"
$ws.Range("E47").Value2 = "  -0.07%  "

# Row 48
$ws.Range("E48").Value2 = "  -0.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "7.61"
$ws.Range("E50").Value2 = "  -0.84%  "

# Row 51
$ws.Range("E51").Value2 = "  -0.45%  "

